$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 1
    4  = -3
    5  = 3
    6  = 6
    7  = -1
    8  = 3
    9  = -1
    10 = -3
    11 = 5
    12 = 5
    14 = 5
    15 = -1
    16 = -3
    17 = -3
    18 = -6
    19 = 4
    21 = 3
    22 = 1
    23 = 2
    24 = 5
    25 = 6
    26 = 6
    27 = 2
    29 = 4
    30 = -2
    31 = 4
    32 = 4
    33 = -1
    34 = 3
    35 = 2
    36 = -2
    37 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
